# "os separators change to nix-like and del unnecessary strings"
# Adds two new card rows (16 and 17) to the "Список" sheet and moves the
# active selection to B17.
#
# New shared strings are introduced in this order so they land at shared
# string indices 27 ("Conflux"), 28 ("Впитать Мощь"), 29 ("Asha's Favor"),
# matching the target workbook - hence D16 is written before C16.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 16: quantity / runame / set (no hyperlink style on the new "set" cell)
$ws.Cells.Item(16, 4).Value = "Conflux"
$ws.Cells.Item(16, 2).Value = 1
$ws.Cells.Item(16, 3).Value = "Впитать Мощь"

# Row 17: enname / quantity / set
$ws.Cells.Item(17, 1).Value = "Asha's Favor"
$ws.Cells.Item(17, 2).Value = 1
$ws.Cells.Item(17, 4).Value = "Conflux"

# Move the selection to B17, as in the updated workbook.
$ws.Range("B17").Select() | Out-Null
